$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Un-Completed" values for rows 6-8 (column C)
$ws.Range("C6").Value = 33
$ws.Range("C7").Value = 33
$ws.Range("C8").Value = 33

# Update the active selection to match the author's final cursor position
$ws.Range("C12").Select()
